$wb = $excel.ActiveWorkbook

# Reorder sheets: review_info becomes the first sheet, hotel_info the second
$wsReview = $wb.Worksheets.Item("review_info")
$wsReview.Move($wb.Worksheets.Item(1))

# Add a new "State" column to hotel_info, between Hotel_Name and City
$wsHotel = $wb.Worksheets.Item("hotel_info")
$wsHotel.Range("C1").EntireColumn.Insert()
$wsHotel.Range("C1").Value = "State"
$wsHotel.Range("C2").Value = "Louisiana"
